$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 295, shifting existing rows 295-408 down to 296-409.
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row 295 with the new weekly price record.
$ws.Cells.Item(295, 1).Value = 4
$ws.Cells.Item(295, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(295, 3).Value = "Los Lagos"
$ws.Cells.Item(295, 4).Value = 44825
$ws.Cells.Item(295, 5).Value = 10
$ws.Cells.Item(295, 6).Value = 100112008
$ws.Cells.Item(295, 7).Value = "Coliflor"
$ws.Cells.Item(295, 8).Value = "Sin especificar"
$ws.Cells.Item(295, 9).Value = "Primera"
$ws.Cells.Item(295, 10).Value = 500
$ws.Cells.Item(295, 11).Value = 1800
$ws.Cells.Item(295, 12).Value = 1800
$ws.Cells.Item(295, 13).Value = 1800
$ws.Cells.Item(295, 14).Value = "$/unidad"
$ws.Cells.Item(295, 15).Value = "Región Metropolitana"
$ws.Cells.Item(295, 16).Value = 1800
$ws.Cells.Item(295, 17).Value = 1
$ws.Cells.Item(295, 18).Value = "Hortaliza"
